$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new data row for the new date (serial 45695 -> 2025-02-07),
# reusing the date style already applied to the row above (A9).
$ws.Range("A10").Value = 45695
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B10").Value = "1 hour"
$ws.Range("C10").Value = "clean properties data"

# Update selection to B2, matching the saved view state
$ws.Range("B2").Select()
